$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: new diary entry for 1/23/2020, 5:00 PM
$ws.Range("A14").Value = 43853
$ws.Range("A14").NumberFormat = "MM/DD/YY"
$ws.Range("B14").Value = 0.708333333333333
$ws.Range("B14").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("C14").Value = "Class"
$ws.Range("D14").Value = "Attend lecture"
$ws.Range("E14").Value = "Learned about mental models and UML diagrams"
$ws.Range("F14").Value = "The mental model is a useful abstraction layer between the model and the code and good for understanding limitations and possible mistakes"
$ws.Range("G14").Value = "Positive"
$ws.Rows.Item(14).RowHeight = 72.35

# Row 15: new diary entry for 1/25/2020, 3:00 PM
$ws.Range("A15").Value = 43855
$ws.Range("A15").NumberFormat = "MM/DD/YY"
$ws.Range("B15").Value = 0.625
$ws.Range("B15").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("C15").Value = "Team"
$ws.Range("D15").Value = "Find 2 features"
$ws.Range("E15").Value = "All goals"
$ws.Range("F15").Value = "UML diagrams can be helpful but can be too complex to use easily. Searching for how a feature is implemented can result in going through a lot of different parts of the program."
$ws.Range("G15").Value = "A lot of stuff going on, but manageable"
$ws.Rows.Item(15).RowHeight = 87

# Move the active selection to G15 as in the updated workbook
$ws.Range("G15").Select()

$wb.Save()
